# Auto-generated edit script applying scheduled-runner market/profit data refresh
# to the Leve profit tracker workbook (sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 674.8570999999999
$ws.Range("I28").Value = 561.6842
$ws.Range("J28").Value = 1750
$ws.Range("K28").Value = 561.6842
$ws.Range("L28").Value = 1750
$ws.Range("M28").Value = -76.68420000000003
$ws.Range("N28").Value = -2720
$ws.Range("H40").Value = 6000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 6000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -6350
$ws.Range("H64").Value = 3563.5117
$ws.Range("J64").Value = 3530.3794
$ws.Range("L64").Value = 3530.3794
$ws.Range("N64").Value = -4026.3794
$ws.Range("H67").Value = 3563.5117
$ws.Range("J67").Value = 3530.3794
$ws.Range("L67").Value = 3530.3794
$ws.Range("N67").Value = -5246.3794
$ws.Range("H70").Value = 2940.5557
$ws.Range("I70").Value = 4338.5557
$ws.Range("J70").Value = 1542.5555
$ws.Range("K70").Value = 13015.6671
$ws.Range("L70").Value = 4627.666499999999
$ws.Range("M70").Value = -12745.6671
$ws.Range("N70").Value = -5167.666499999999
$ws.Range("H73").Value = 2940.5557
$ws.Range("I73").Value = 4338.5557
$ws.Range("J73").Value = 1542.5555
$ws.Range("K73").Value = 13015.6671
$ws.Range("L73").Value = 4627.666499999999
$ws.Range("M73").Value = -12079.6671
$ws.Range("N73").Value = -6499.666499999999
$ws.Range("H74").Value = 4138.55
$ws.Range("I74").Value = 4277.5713
$ws.Range("J74").Value = 4063.6924
$ws.Range("K74").Value = 4277.5713
$ws.Range("L74").Value = 4063.6924
$ws.Range("M74").Value = -3341.5713
$ws.Range("N74").Value = -5935.6924
$ws.Range("H76").Value = 8957.458000000001
$ws.Range("I76").Value = 12256.583
$ws.Range("J76").Value = 5658.3335
$ws.Range("K76").Value = 12256.583
$ws.Range("L76").Value = 5658.3335
$ws.Range("M76").Value = -11941.583
$ws.Range("N76").Value = -6288.3335
$ws.Range("H77").Value = 4138.55
$ws.Range("I77").Value = 4277.5713
$ws.Range("J77").Value = 4063.6924
$ws.Range("K77").Value = 21387.8565
$ws.Range("L77").Value = 20318.462
$ws.Range("M77").Value = -16707.8565
$ws.Range("N77").Value = -29678.462
$ws.Range("H79").Value = 8957.458000000001
$ws.Range("I79").Value = 12256.583
$ws.Range("J79").Value = 5658.3335
$ws.Range("K79").Value = 12256.583
$ws.Range("L79").Value = 5658.3335
$ws.Range("M79").Value = -11164.583
$ws.Range("N79").Value = -7842.3335
$ws.Range("H98").Value = 1089.9615
$ws.Range("I98").Value = 1025.381
$ws.Range("J98").Value = 1361.2
$ws.Range("K98").Value = 1025.381
$ws.Range("L98").Value = 1361.2
$ws.Range("M98").Value = 472.6189999999999
$ws.Range("N98").Value = -4357.2
$ws.Range("H111").Value = 1019.93335
$ws.Range("I111").Value = 785.2222
$ws.Range("J111").Value = 1372
$ws.Range("K111").Value = 2355.6666
$ws.Range("L111").Value = 4116
$ws.Range("M111").Value = 711.3334
$ws.Range("N111").Value = -10250
$ws.Range("H122").Value = 1089.9615
$ws.Range("I122").Value = 1025.381
$ws.Range("J122").Value = 1361.2
$ws.Range("K122").Value = 3076.143
$ws.Range("L122").Value = 4083.6
$ws.Range("M122").Value = -626.143
$ws.Range("N122").Value = -8983.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 441
$ws.Range("I2").Value = 355.66666
$ws.Range("J2").Value = 825
$ws.Range("K2").Value = 355.66666
$ws.Range("L2").Value = 825
$ws.Range("M2").Value = -242.66666
$ws.Range("N2").Value = -1051
$ws.Range("H63").Value = 4958.407
$ws.Range("I63").Value = 5578.9
$ws.Range("J63").Value = 3185.5715
$ws.Range("K63").Value = 5578.9
$ws.Range("L63").Value = 3185.5715
$ws.Range("M63").Value = -4892.9
$ws.Range("N63").Value = -4557.5715
$ws.Range("H66").Value = 4958.407
$ws.Range("I66").Value = 5578.9
$ws.Range("J66").Value = 3185.5715
$ws.Range("K66").Value = 27894.5
$ws.Range("L66").Value = 15927.8575
$ws.Range("M66").Value = -24462.5
$ws.Range("N66").Value = -22791.8575
$ws.Range("H116").Value = 441
$ws.Range("I116").Value = 355.66666
$ws.Range("J116").Value = 825
$ws.Range("K116").Value = 355.66666
$ws.Range("L116").Value = 825
$ws.Range("M116").Value = 1938.33334
$ws.Range("N116").Value = -5413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 441
$ws.Range("I3").Value = 355.66666
$ws.Range("J3").Value = 825
$ws.Range("K3").Value = 355.66666
$ws.Range("L3").Value = 825
$ws.Range("M3").Value = -241.66666
$ws.Range("N3").Value = -1053
$ws.Range("H94").Value = 893.2281
$ws.Range("I94").Value = 795.53656
$ws.Range("J94").Value = 1143.5625
$ws.Range("K94").Value = 795.53656
$ws.Range("L94").Value = 1143.5625
$ws.Range("M94").Value = -344.53656
$ws.Range("N94").Value = -2045.5625
$ws.Range("H105").Value = 2395.2144
$ws.Range("I105").Value = 2309.5454
$ws.Range("K105").Value = 2309.5454
$ws.Range("M105").Value = -562.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1404.7544
$ws.Range("I134").Value = 809.2683
$ws.Range("K134").Value = 2427.8049
$ws.Range("M134").Value = 107.1950999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 351.61905
$ws.Range("J5").Value = 1421
$ws.Range("K5").Value = 1054.85715
$ws.Range("L5").Value = 4263
$ws.Range("M5").Value = -942.85715
$ws.Range("N5").Value = -4487
$ws.Range("H23").Value = 506.32
$ws.Range("J23").Value = 525.7083
$ws.Range("L23").Value = 1577.1249
$ws.Range("N23").Value = -2047.1249
$ws.Range("I135").Value = 351.61905
$ws.Range("J135").Value = 1421
$ws.Range("K135").Value = 3164.57145
$ws.Range("L135").Value = 12789
$ws.Range("M135").Value = -629.5714500000004
$ws.Range("N135").Value = -17859

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10475.389
$ws.Range("I70").Value = 3896.4167
$ws.Range("K70").Value = 3896.4167
$ws.Range("M70").Value = -3626.4167
$ws.Range("H73").Value = 10475.389
$ws.Range("I73").Value = 3896.4167
$ws.Range("K73").Value = 3896.4167
$ws.Range("M73").Value = -2960.4167
$ws.Range("H80").Value = 4318.3257
$ws.Range("I80").Value = 4647.0586
$ws.Range("J80").Value = 3076.4443
$ws.Range("K80").Value = 4647.0586
$ws.Range("L80").Value = 3076.4443
$ws.Range("M80").Value = -3649.0586
$ws.Range("N80").Value = -5072.4443
$ws.Range("H83").Value = 4318.3257
$ws.Range("I83").Value = 4647.0586
$ws.Range("J83").Value = 3076.4443
$ws.Range("K83").Value = 23235.293
$ws.Range("L83").Value = 15382.2215
$ws.Range("M83").Value = -18243.293
$ws.Range("N83").Value = -25366.2215
$ws.Range("H102").Value = 3372.9048
$ws.Range("I102").Value = 3753.4
$ws.Range("J102").Value = 1470.4286
$ws.Range("K102").Value = 3753.4
$ws.Range("L102").Value = 1470.4286
$ws.Range("M102").Value = -2131.4
$ws.Range("N102").Value = -4714.4286
$ws.Range("H107").Value = 5211.6
$ws.Range("I107").Value = 7276.5
$ws.Range("J107").Value = 393.5
$ws.Range("K107").Value = 7276.5
$ws.Range("L107").Value = 393.5
$ws.Range("M107").Value = -5356.5
$ws.Range("N107").Value = -4233.5
$ws.Range("H122").Value = 1733.2222
$ws.Range("I122").Value = 1683.1666
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 5049.4998
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -2599.4998
$ws.Range("N122").Value = -10400.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2854.8333
$ws.Range("I40").Value = 2630.8
$ws.Range("K40").Value = 2630.8
$ws.Range("M40").Value = -2494.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 615.4400000000001
$ws.Range("I107").Value = 739.4
$ws.Range("J107").Value = 532.8
$ws.Range("K107").Value = 2218.2
$ws.Range("L107").Value = 1598.4
$ws.Range("M107").Value = -298.1999999999998
$ws.Range("N107").Value = -5438.4
